# Jogos_da_Semana_FlashScore_2024-11-16.xlsx update
# 1) Remove the Paraguay "Cerro Porteno - Guarani" match (old row 4); this shifts
#    the two USL Championship rows up by one (old rows 5,6 -> new rows 4,5).
# 2) Refresh the betting odds (columns G:BD) for rows 2, 3 and the new row 4
#    (the "Charleston - Rhode Island" match) to the latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete the Paraguay match row ---
$ws.Rows.Item(4).Delete()

# --- Step 2: refresh odds values ---
$cols = @("G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD")

$row2Vals = @(1.73,3.6,5,2.38,2.2,5,1.05,11,1.29,3.5,1.95,1.9,1.4,2.75,1.8,1.91,7.5,8.5,8.5,13,13,26,10,7,15,51,251,13,26,15,51,41,41,3.75,9,21,29,51,151,2.75,8,51,6.5,23,34,81,101,201,81,81)
$row3Vals = @(4.5,3.2,1.9,5,2,2.63,1.1,7,1.4,2.75,2.35,1.57,1.5,2.5,2.1,1.67,10,21,15,51,41,51,7,6.5,19,67,351,6,8,9,15,17,34,6,26,41,101,126,351,2.5,9,67,3.75,11,26,41,67,201,81,81)
$row4Vals = @(2.12,3.2,3.15,2.72,2.1,3.7,1.05,7.7,1.25,3.6,1.75,2,1.38,2.8,1.6,2.2,8.75,11.5,8.5,22,16,22,7.7,6.4,11.75,45,300,11.25,18.5,11,45,26,29,4.2,11.25,18,45,70,200,2.8,6.6,50,5.3,17.5,22,90,110,250,51,51)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $row2Vals[$i]
    $ws.Range($cols[$i] + "3").Value = $row3Vals[$i]
    $ws.Range($cols[$i] + "4").Value = $row4Vals[$i]
}
